$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 3549.625
$ws.Cells.Item(40, 9).Value = 2867
$ws.Cells.Item(40, 11).Value = 2867
$ws.Cells.Item(40, 13).Value = -2692
# Row 64
$ws.Cells.Item(64, 8).Value = 4860.25
$ws.Cells.Item(64, 9).Value = 4860.25
$ws.Cells.Item(64, 11).Value = 4860.25
$ws.Cells.Item(64, 13).Value = -4612.25
# Row 67
$ws.Cells.Item(67, 8).Value = 4860.25
$ws.Cells.Item(67, 9).Value = 4860.25
$ws.Cells.Item(67, 11).Value = 4860.25
$ws.Cells.Item(67, 13).Value = -4002.25
# Row 70
$ws.Cells.Item(70, 8).Value = 11333.333
$ws.Cells.Item(70, 10).Value = 12000
$ws.Cells.Item(70, 12).Value = 36000
$ws.Cells.Item(70, 14).Value = -36540
# Row 73
$ws.Cells.Item(73, 8).Value = 11333.333
$ws.Cells.Item(73, 10).Value = 12000
$ws.Cells.Item(73, 12).Value = 36000
$ws.Cells.Item(73, 14).Value = -37872
# Row 100
$ws.Cells.Item(100, 8).Value = 2915.6667
$ws.Cells.Item(100, 9).Value = 2915.6667
$ws.Cells.Item(100, 11).Value = 2915.6667
$ws.Cells.Item(100, 13).Value = -2374.6667
# Row 113
$ws.Cells.Item(113, 8).Value = 4066.4666
$ws.Cells.Item(113, 9).Value = 4178.4287
$ws.Cells.Item(113, 11).Value = 4178.4287
$ws.Cells.Item(113, 13).Value = -924.4287000000004
# Row 127
$ws.Cells.Item(127, 8).Value = 1726.6
$ws.Cells.Item(127, 9).Value = 1726.6
$ws.Cells.Item(127, 11).Value = 5179.799999999999
$ws.Cells.Item(127, 13).Value = -219.7999999999993
# Row 138
$ws.Cells.Item(138, 8).Value = 1358.625
$ws.Cells.Item(138, 9).Value = 1358.625
$ws.Cells.Item(138, 11).Value = 4075.875
$ws.Cells.Item(138, 13).Value = 1064.125
# Row 141
$ws.Cells.Item(141, 8).Value = 5311.65
$ws.Cells.Item(141, 9).Value = 5796.278
$ws.Cells.Item(141, 10).Value = 950
$ws.Cells.Item(141, 11).Value = 17388.834
$ws.Cells.Item(141, 12).Value = 2850
$ws.Cells.Item(141, 13).Value = -12208.834
$ws.Cells.Item(141, 14).Value = -13210

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2866.3333
$ws.Cells.Item(45, 10).Value = 2000
$ws.Cells.Item(45, 12).Value = 2000
$ws.Cells.Item(45, 14).Value = -2754
# Row 132
$ws.Cells.Item(132, 8).Value = 4301.3105
$ws.Cells.Item(132, 9).Value = 4311.64
$ws.Cells.Item(132, 10).Value = 4236.75
$ws.Cells.Item(132, 11).Value = 12934.92
$ws.Cells.Item(132, 12).Value = 12710.25
$ws.Cells.Item(132, 13).Value = -10404.92
$ws.Cells.Item(132, 14).Value = -17770.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 1112.6
$ws.Cells.Item(107, 9).Value = 666
$ws.Cells.Item(107, 11).Value = 666
$ws.Cells.Item(107, 13).Value = 1254

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1692.125
$ws.Cells.Item(31, 9).Value = 1682.409
$ws.Cells.Item(31, 10).Value = 1799
$ws.Cells.Item(31, 11).Value = 1682.409
$ws.Cells.Item(31, 12).Value = 1799
$ws.Cells.Item(31, 13).Value = -1387.409
$ws.Cells.Item(31, 14).Value = -2389
# Row 34
$ws.Cells.Item(34, 8).Value = 1692.125
$ws.Cells.Item(34, 9).Value = 1682.409
$ws.Cells.Item(34, 10).Value = 1799
$ws.Cells.Item(34, 11).Value = 1682.409
$ws.Cells.Item(34, 12).Value = 1799
$ws.Cells.Item(34, 13).Value = -1480.409
$ws.Cells.Item(34, 14).Value = -2203
# Row 36
$ws.Cells.Item(36, 8).Value = 4000
$ws.Cells.Item(36, 9).Value = 4000
$ws.Cells.Item(36, 11).Value = 4000
$ws.Cells.Item(36, 13).Value = -3612
# Row 40
$ws.Cells.Item(40, 8).Value = 4000
$ws.Cells.Item(40, 9).Value = 4000
$ws.Cells.Item(40, 11).Value = 4000
$ws.Cells.Item(40, 13).Value = -3840
# Row 58
$ws.Cells.Item(58, 8).Value = 2114.5881
$ws.Cells.Item(58, 9).Value = 1564.1428
$ws.Cells.Item(58, 10).Value = 4683.3335
$ws.Cells.Item(58, 11).Value = 1564.1428
$ws.Cells.Item(58, 12).Value = 4683.3335
$ws.Cells.Item(58, 13).Value = -1361.1428
$ws.Cells.Item(58, 14).Value = -5089.3335
# Row 59
$ws.Cells.Item(59, 8).Value = 83597.89
$ws.Cells.Item(59, 10).Value = 85297.625
$ws.Cells.Item(59, 12).Value = 85297.625
$ws.Cells.Item(59, 14).Value = -87587.625
# Row 62
$ws.Cells.Item(62, 8).Value = 10999.5
$ws.Cells.Item(62, 9).Value = 10999.5
$ws.Cells.Item(62, 11).Value = 10999.5
$ws.Cells.Item(62, 13).Value = -10375.5
# Row 65
$ws.Cells.Item(65, 8).Value = 10999.5
$ws.Cells.Item(65, 9).Value = 10999.5
$ws.Cells.Item(65, 11).Value = 54997.5
$ws.Cells.Item(65, 13).Value = -51877.5
# Row 94
$ws.Cells.Item(94, 8).Value = 1609.8
$ws.Cells.Item(94, 9).Value = 1687.25
$ws.Cells.Item(94, 10).Value = 1300
$ws.Cells.Item(94, 11).Value = 1687.25
$ws.Cells.Item(94, 12).Value = 1300
$ws.Cells.Item(94, 13).Value = -1236.25
$ws.Cells.Item(94, 14).Value = -2202
# Row 136
$ws.Cells.Item(136, 8).Value = 2114.5881
$ws.Cells.Item(136, 9).Value = 1564.1428
$ws.Cells.Item(136, 10).Value = 4683.3335
$ws.Cells.Item(136, 11).Value = 4692.428400000001
$ws.Cells.Item(136, 12).Value = 14050.0005
$ws.Cells.Item(136, 13).Value = -2142.428400000001
$ws.Cells.Item(136, 14).Value = -19150.0005

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 846.35
$ws.Cells.Item(5, 10).Value = 848.1429000000001
$ws.Cells.Item(5, 12).Value = 2544.4287
$ws.Cells.Item(5, 14).Value = -2768.4287
# Row 131
$ws.Cells.Item(131, 8).Value = 1786.4
$ws.Cells.Item(131, 10).Value = 2160.6667
$ws.Cells.Item(131, 12).Value = 6482.000100000001
$ws.Cells.Item(131, 14).Value = -16562.0001
# Row 135
$ws.Cells.Item(135, 8).Value = 846.35
$ws.Cells.Item(135, 10).Value = 848.1429000000001
$ws.Cells.Item(135, 12).Value = 7633.2861
$ws.Cells.Item(135, 14).Value = -12703.2861

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 3446.7856
$ws.Cells.Item(122, 9).Value = 3014.0908
$ws.Cells.Item(122, 10).Value = 5033.3335
$ws.Cells.Item(122, 11).Value = 9042.2724
$ws.Cells.Item(122, 12).Value = 15100.0005
$ws.Cells.Item(122, 13).Value = -6592.2724
$ws.Cells.Item(122, 14).Value = -20000.0005

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 2316.6
$ws.Cells.Item(22, 9).Value = 2395.75
$ws.Cells.Item(22, 11).Value = 2395.75
$ws.Cells.Item(22, 13).Value = -2100.75
# Row 27
$ws.Cells.Item(27, 8).Value = 2316.6
$ws.Cells.Item(27, 9).Value = 2395.75
$ws.Cells.Item(27, 11).Value = 2395.75
$ws.Cells.Item(27, 13).Value = -2288.75
# Row 40
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).ClearContents()
# Row 76
$ws.Cells.Item(76, 8).Value = 15999
$ws.Cells.Item(76, 10).Value = 15999
$ws.Cells.Item(76, 12).Value = 15999
$ws.Cells.Item(76, 14).Value = -16675
# Row 79
$ws.Cells.Item(79, 8).Value = 15999
$ws.Cells.Item(79, 10).Value = 15999
$ws.Cells.Item(79, 12).Value = 15999
$ws.Cells.Item(79, 14).Value = -18339
# Row 122
$ws.Cells.Item(122, 8).Value = 3121.75
$ws.Cells.Item(122, 9).Value = 2995
$ws.Cells.Item(122, 10).Value = 3248.5
$ws.Cells.Item(122, 11).Value = 8985
$ws.Cells.Item(122, 12).Value = 9745.5
$ws.Cells.Item(122, 13).Value = -6535
$ws.Cells.Item(122, 14).Value = -14645.5
# Row 132
$ws.Cells.Item(132, 8).Value = 2364.6365
$ws.Cells.Item(132, 9).Value = 2365.2727
$ws.Cells.Item(132, 10).Value = 2363.3635
$ws.Cells.Item(132, 11).Value = 7095.8181
$ws.Cells.Item(132, 12).Value = 7090.0905
$ws.Cells.Item(132, 13).Value = -4565.8181
$ws.Cells.Item(132, 14).Value = -12150.0905

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Cells.Item(45, 8).Value = 42555.363
$ws.Cells.Item(45, 9).Value = 31674
$ws.Cells.Item(45, 10).Value = 55613
$ws.Cells.Item(45, 11).Value = 31674
$ws.Cells.Item(45, 12).Value = 55613
$ws.Cells.Item(45, 13).Value = -31183
$ws.Cells.Item(45, 14).Value = -56595
# Row 74
$ws.Cells.Item(74, 8).Value = 19587
$ws.Cells.Item(74, 9).Value = 18645
$ws.Cells.Item(74, 11).Value = 18645
$ws.Cells.Item(74, 13).Value = -17709
# Row 77
$ws.Cells.Item(77, 8).Value = 19587
$ws.Cells.Item(77, 9).Value = 18645
$ws.Cells.Item(77, 11).Value = 55935
$ws.Cells.Item(77, 13).Value = -51255
# Row 126
$ws.Cells.Item(126, 8).Value = 1854.1111
$ws.Cells.Item(126, 9).Value = 1210.875
$ws.Cells.Item(126, 11).Value = 3632.625
$ws.Cells.Item(126, 13).Value = -1162.625
# Row 132
$ws.Cells.Item(132, 8).Value = 18791.6
$ws.Cells.Item(132, 9).Value = 18334.666
$ws.Cells.Item(132, 11).Value = 55003.99800000001
$ws.Cells.Item(132, 13).Value = -52473.99800000001
# Row 136
$ws.Cells.Item(136, 8).Value = 8447.909
$ws.Cells.Item(136, 9).Value = 8447.909
$ws.Cells.Item(136, 11).Value = 25343.727
$ws.Cells.Item(136, 13).Value = -22793.727
